$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values (test17 -> test14 environment, new creds) ---

# A2 (BOURL): backoffice URL test17 -> test14
$ws.Range("A2").Value = "https://test14.cliotest.com/backoffice/control/main"

# E2 (ofbizuser): abcd -> sshinde
$ws.Range("E2").Value = "sshinde"

# F2 (ofbizpass): was numeric 1234, now the text password (keeps style s=3)
$ws.Range("F2").Value = 'C@bi$ush5'

# G2 (CCURL): cabicentral URL test17 -> test14
$ws.Range("G2").Value = "https://test14.cliotest.com/cabicentral/control/main"

# J2 (CWURL): warehouse URL test17 -> test14
$ws.Range("J2").Value = "https://test14.cliotest.com/warehouse/control/main"

# --- Update hyperlink display (TextToDisplay) attributes ---
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$F$2') {
        # mailto hyperlink no longer carries an explicit display override
        $h.TextToDisplay = ""
    } elseif ($addr -eq '$G$2') {
        $h.TextToDisplay = "https://test17.cliotest.com/cabicentral/control/main"
    } elseif ($addr -eq '$J$2') {
        $h.TextToDisplay = "https://test17.cliotest.com/warehouse/control/main"
    } elseif ($addr -eq '$A$2') {
        $h.TextToDisplay = "https://test17.cliotest.com/backoffice/control/main"
    }
}

# --- Update selection / view state: F2 selected, no pinned top-left cell ---
$ws.Range("F2").Select()
